# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) held a mangled literal string ("5-12-2013-14")
# that should read "2014-05-12" for every data row (rows 2-31).
#
# Note: assigning the literal text "2014-05-12" straight to .Value would
# get auto-recognized as a real Excel date serial (since it parses as a
# valid ISO date) instead of staying a plain text string like the source
# file expects. Building it via a text formula and then collapsing the
# formula down to its static result with PasteSpecial(values) keeps it a
# literal string, matching the original "Date" column formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "5-12-2013-14"
$newDate = "2014-05-12"

$dateRange = $ws.Range("BF2:BF31")

$allMatch = $true
for ($row = 2; $row -le 31; $row++) {
    if ($ws.Range("BF$row").Value() -ne $oldDate) {
        $allMatch = $false
    }
}

if ($allMatch) {
    $dateRange.Formula = '="' + $newDate + '"'
    $dateRange.Copy()
    $dateRange.PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = 0
}
